# Regenerate merged AHB files
#
# Rows 93-126 in the sole worksheet belong to the same "merged AHB" table
# block. This edit:
#   1) Re-styles the group-header rows (the first row of each item group,
#      where column A carries the new item number) from the plain bordered
#      style (s=5 / s=5 / s=5) to the shaded-row style (s=2 / s=3(col B) /
#      s=2), matching the style already used e.g. by row 2.
#   2) Clears the "AENDERUNG" marker text out of column L for every row in
#      93-126 and re-styles those L cells from s=7 (bold, gold, centered)
#      to s=4 (centered, same shading, default font) -- again matching the
#      style already used e.g. by L2.
#
# We drive this purely through Copy / PasteSpecial(xlPasteFormats) from
# cells that already carry the exact target style, so no new style (xf)
# records get minted -- the engine dedupes onto the existing cellXfs
# entries (2, 3, 4) instead of creating near-duplicates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Group-header rows: first row of each item group in the 93-126 block.
$headerRows = @(93, 96, 100, 103, 107, 111, 115, 121, 124)

# Full range of rows touched in this block (header rows + the detail rows
# that only lose their column-L marker).
$firstRow = 93
$lastRow = 126

# --- Step 1: re-style the group-header rows (columns A:V) -------------
# Row 2 already has exactly the target look (A:V = s2/s3/s2...), so copy
# its formatting once and paste it onto each header row.
$ws.Range("A2:V2").Copy() | Out-Null
foreach ($r in $headerRows) {
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Step 2: re-style column L for every row in the block --------------
# L2 already has exactly the target look (s=4, blank), so copy its
# formatting once and paste across the whole contiguous L93:L126 range.
$ws.Range("L2").Copy() | Out-Null
$ws.Range("L" + $firstRow + ":L" + $lastRow).PasteSpecial($xlPasteFormats) | Out-Null

# --- Step 3: clear the "AENDERUNG" text that used to live in column L --
$ws.Range("L" + $firstRow + ":L" + $lastRow).ClearContents() | Out-Null

$excel.CutCopyMode = $false
